$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45203 -> 45204) for every data row (rows 2 through 374).
$ws.Range("C2:C374").Value = 45204
